$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (D and G hold numeric-looking text values)
# so Excel does not auto-convert them to numbers on write.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '245.42'
$ws.Range("G2").Value = '9'

# Row 3
$ws.Range("D3").Value = '25.92'
$ws.Range("G3").Value = '9'

# Row 4
$ws.Range("D4").Value = '5.123'
$ws.Range("G4").Value = '9'

# Row 5
$ws.Range("D5").Value = '0.05592'
$ws.Range("G5").Value = '9'

# Row 6
$ws.Range("D6").Value = '6.482'
$ws.Range("G6").Value = '9'

# Row 7
$ws.Range("D7").Value = '3.029'
$ws.Range("G7").Value = '9'

# Row 8
$ws.Range("D8").Value = '0.8162'
$ws.Range("G8").Value = '9'

# Row 9
$ws.Range("D9").Value = '0.8464'
$ws.Range("G9").Value = '9'

# Row 10
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '0.009713'
$ws.Range("E10").Value = '9OneONEBestin24h'
$ws.Range("G10").Value = '9'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1344'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("G11").Value = '9'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.02852'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("G12").Value = '9'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09408'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("G13").Value = '9'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001514'
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("G14").Value = '9'

# Row 15
$ws.Range("D15").Value = '0.006131'
$ws.Range("G15").Value = '9'

# Row 16
$ws.Range("D16").Value = '3.550'
$ws.Range("G16").Value = '9'

# Row 17
$ws.Range("G17").Value = '9'

# Row 18
$ws.Range("G18").Value = '9'

# Row 19
$ws.Range("D19").Value = '0.06953'
$ws.Range("G19").Value = '9'

# Row 20
$ws.Range("D20").Value = '0.03234'
$ws.Range("G20").Value = '9'

# Row 21
$ws.Range("G21").Value = '9'

# Row 22
$ws.Range("D22").Value = '3.754'
$ws.Range("G22").Value = '9'

# Row 23
$ws.Range("D23").Value = '0.04699'
$ws.Range("G23").Value = '9'

# Row 24
$ws.Range("G24").Value = '9'

# Row 25
$ws.Range("D25").Value = '0.001251'
$ws.Range("G25").Value = '9'

# Row 26
$ws.Range("D26").Value = '0.004604'
$ws.Range("G26").Value = '9'

# Row 27
$ws.Range("D27").Value = '0.00009605'
$ws.Range("G27").Value = '9'

# Row 28
$ws.Range("D28").Value = '0.0001391'
$ws.Range("G28").Value = '9'

# Row 29
$ws.Range("G29").Value = '9'

# Row 30
$ws.Range("G30").Value = '9'

# Row 31
$ws.Range("G31").Value = '9'

# Row 32
$ws.Range("G32").Value = '9'

# Row 33
$ws.Range("G33").Value = '9'

# Row 34
$ws.Range("G34").Value = '9'

# Row 35
$ws.Range("G35").Value = '9'

# Row 36
$ws.Range("G36").Value = '9'

# Row 37
$ws.Range("G37").Value = '9'

# Row 38
$ws.Range("G38").Value = '9'

# Row 39
$ws.Range("G39").Value = '9'

# Row 40
$ws.Range("D40").Value = '0.03657'
$ws.Range("G40").Value = '9'

# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '0.006113'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("G41").Value = '9'

# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1053'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("G42").Value = '9'

# Row 43
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '0.002268'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("G43").Value = '9'

# Row 44
$ws.Range("D44").Value = '0.007897'
$ws.Range("G44").Value = '9'

# Row 45
$ws.Range("D45").Value = '0.00005321'
$ws.Range("G45").Value = '9'

# Row 46
$ws.Range("G46").Value = '9'

# Row 47
$ws.Range("D47").Value = '0.1336'
$ws.Range("G47").Value = '9'

# Row 48
$ws.Range("D48").Value = '0.002129'
$ws.Range("G48").Value = '9'

# Row 49
$ws.Range("D49").Value = '0.00002101'
$ws.Range("G49").Value = '9'

# Row 50
$ws.Range("D50").Value = '0.0002001'
$ws.Range("G50").Value = '9'

# Row 51
$ws.Range("G51").Value = '9'
